# Bug fix: palette panel overlapping image and status panel
#
# The "Active" (Todo) sheet had a bug entry (Id 42) for the palette panel
# overlapping the picturebox/statuspanel. That bug is now fixed, so its row
# is removed from "Active" and a corresponding "Done" row is inserted at the
# top of the "Inactive" sheet. In its place, a brand new follow-up Todo task
# (building a MasterImage class) is appended to "Active" with the next Id
# (54), and the "Config" sheet's running Max Id counter is bumped to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Active" sheet: drop the now-fixed bug row (row 3, Id 42) ...
# ---------------------------------------------------------------------
$active = $wb.Worksheets.Item("Active")
$active.Rows.Item(3).Delete()

# ...and append a new Todo task at the end of the list (row 9, right before
# the existing Id-53 row) describing the MasterImage work.
$active.Rows.Item(9).Insert()

$newTaskText = "make a MasterImage class`nthat one-time figures out all the regions in a background thread`nand provides that data to the color worker"

$active.Cells.Item(9, 1).Value = 54
$active.Cells.Item(9, 2).Value = $newTaskText
$active.Cells.Item(9, 3).Value = "Todo"
$active.Cells.Item(9, 4).Value = "Task"
$active.Cells.Item(9, 5).NumberFormat = "@"
$active.Cells.Item(9, 5).Value = "8/23/2018"

# ---------------------------------------------------------------------
# 2) "Inactive" sheet: insert the fixed bug at the top as a Done item.
# ---------------------------------------------------------------------
$inactive = $wb.Worksheets.Item("Inactive")
$inactive.Rows.Item(2).Insert()

$inactive.Cells.Item(2, 1).Value = 42
$inactive.Cells.Item(2, 2).Value = "bug: expanded palette covers part of picturebox and statuspanel"
$inactive.Cells.Item(2, 3).Value = "Done"
$inactive.Cells.Item(2, 4).Value = "Task"
$inactive.Cells.Item(2, 5).NumberFormat = "@"
$inactive.Cells.Item(2, 5).Value = "8/15/2018"
$inactive.Cells.Item(2, 6).NumberFormat = "@"
$inactive.Cells.Item(2, 6).Value = "8/23/2018"

# ---------------------------------------------------------------------
# 3) "Config" sheet: bump the Max Id counter for Active/Todo/Task to 54.
# ---------------------------------------------------------------------
$config = $wb.Worksheets.Item("Config")
$config.Cells.Item(2, 6).Value = 54
